$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New sval data (filtered save games) for rows 2-9, columns B:E and G (sum)
$data = @{
    2 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    3 = @{ B = 3.230985683306322;  C = 10.29869402782916;  D = 0.1575252929769615; E = 8.660232485948974;  G = 22.34743749006142 }
    4 = @{ B = 0.6753301551942219; C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 3.645393585217082 }
    5 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 0.8054896365839992; E = 0.496779210170732;  G = 6.201049113329182 }
    6 = @{ B = 0.01514828764759746; C = 0.3127903958511391; D = 0.1575252929769615; E = 0.496779210170732;  G = 0.9822431866464301 }
    7 = @{ B = 1.459612070389937;  C = 1.667794583268128;  D = 0.1575252929769615; E = 8.660232485948974;  G = 11.945164432584 }
    8 = @{ B = 0.04763786555579896; C = 0.3127903958511391; D = 0.8054896365839992; E = 8.660232485948974;  G = 9.826150383939911 }
    9 = @{ B = 3.230985683306322;  C = 1.667794583268128;  D = 337.1190423067083;  E = 8.660232485948974;  G = 350.6780550592317 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 2).Value = $vals.B
    $ws.Cells.Item($row, 3).Value = $vals.C
    $ws.Cells.Item($row, 4).Value = $vals.D
    $ws.Cells.Item($row, 5).Value = $vals.E
    $ws.Cells.Item($row, 7).Value = $vals.G
}
